$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update account number on first data row (Excel row 2, LUIS)
#    Keep it as text so the leading zeros are preserved.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "008007764"

# 2. Delete the entire row for MARIA / 004204242 / 30110.62 (Excel row 4)
$ws.Rows.Item(4).Delete()

# The sheet has a completely empty placeholder row further down (originally
# Excel row 384, now shifted to row 383 after the delete above). The engine's
# row-shift does not re-materialize rows that carry no cells at all, so nudge
# a row-level (not cell-level) property to force it back into existence
# without altering its appearance.
$ws.Rows.Item(383).OutlineLevel = 0

# 3. Update balance for BEATRIZ (005046790) from 10000 to 9000.
#    After the row deletion above, this row shifted up from row 7 to row 6.
$ws.Range("C6").Value = 9000
